$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content changes -------------------------------------------------
# G22: "SEGURADORAS" -> "SEGURADORA"
$ws.Range("G22").Value = "SEGURADORA"

# G17: new cell "DTA_FIM" (previously empty), matching the thin-border
# style ("s=1") used by the rest of the table's data cells.
$ws.Range("G17").Value = "DTA_FIM"
$ws.Range("G17").Borders.LineStyle = 1
$ws.Range("G17").Borders.Weight = 2

# I17: "NUM_CPF_FUNCIONARIO <FK> <PK>" -> "NUM_CPF_FUNCIONARIO <FK> "
$ws.Range("I17").Value = "NUM_CPF_FUNCIONARIO <FK> "

# --- Header font color: red (FF0000) -> purple (7030A0) --------------
$headerCells = @("A1","C1","E1","G1","I1","A13","C13","E13","G13","I13","A22","C22","E22","G22","I22")
foreach ($cell in $headerCells) {
    $ws.Range($cell).Font.Color = 10498160
}

# --- Selection moves from G9 to J15 -----------------------------------
$ws.Range("J15").Select()
